$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(78, 1).Value = "2024-10-19 00:00:00"
$ws.Cells.Item(78, 2).Value = 73650
$ws.Cells.Item(78, 3).Value = 10320.18
$ws.Cells.Item(78, 4).Value = 9132.91
$ws.Cells.Item(78, 5).Value = 7.1018
